$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 9968.7199999999993
$ws.Range("B5").Value = 10031.92
$ws.Range("C5").Value = 307.87
$ws.Range("D5").Value = 305.93
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -0.63
$ws.Range("G5").Value = 42612.675347222219
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = $false
